# Refresh scraped event data (想去人数 / 最低票价) on both the "展览" sheet
# and the "全部类型" sheet, which both list the same exhibitions (just
# interleaved with different other rows, so the row numbers differ for
# the last two events).

$wb = $excel.ActiveWorkbook

$sheets = @(
    @{ Name = "展览";     SoldOutRow = 2; Row648 = 6; Row1630 = 7 },
    @{ Name = "全部类型"; SoldOutRow = 2; Row648 = 8; Row1630 = 9 }
)

foreach ($info in $sheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）: now sold out
    $ws.Cells.Item($info.SoldOutRow, 6).Value = 16194
    $ws.Cells.Item($info.SoldOutRow, 7).Value = "已售罄"

    # 南宁·火影忍者only
    $ws.Range("F3").Value = 344

    # 南宁·蔚蓝档案only
    $ws.Range("F4").Value = 721

    # 南宁·熊喵M动漫嘉年华【免费】
    $ws.Cells.Item($info.Row648, 6).Value = 666

    # 南宁·第二届北极光动漫展
    $ws.Cells.Item($info.Row1630, 6).Value = 1650
}
